$wb = $excel.ActiveWorkbook

$wsDecisionVars = $wb.Worksheets.Item("optimal_decision_variables")
$wsDecisionVars.Range("C2").Value = 9981.625679984112
$wsDecisionVars.Range("C3").Value = 462.470649717059
$wsDecisionVars.Range("C4").Value = 8944.504220828099
$wsDecisionVars.Range("C5").Value = 4000
$wsDecisionVars.Range("C6").Value = 460.4583777822451
$wsDecisionVars.Range("C7").Value = 4000
$wsDecisionVars.Range("C8").Value = 421.6160198117984
$wsDecisionVars.Range("C9").Value = 3513.984522373709
$wsDecisionVars.Range("C10").Value = 455.7833425069027
$wsDecisionVars.Range("C11").Value = 4000
$wsDecisionVars.Range("C12").Value = 429.6731768325639
$wsDecisionVars.Range("C13").Value = 3687.452875457066
$wsDecisionVars.Range("C14").Value = 460.4583777822451
$wsDecisionVars.Range("C16").Value = 6.01191402680718
$wsDecisionVars.Range("C17").Value = 0.3270748212356701

$wsEffluent = $wb.Worksheets.Item("optimal_predicted_effluent")
$wsEffluent.Range("B2").Value = 607.265214337428
$wsEffluent.Range("B3").Value = 6731.511399869436
$wsEffluent.Range("B4").Value = 208.42783908035
$wsEffluent.Range("B5").Value = 295.8378432015441
$wsEffluent.Range("B6").Value = 295.8378432015442
$wsEffluent.Range("B7").Value = 5917.262147825887
$wsEffluent.Range("B8").Value = 4898.246857605723
$wsEffluent.Range("B9").Value = 60.3273602907266
$wsEffluent.Range("B10").Value = 740.6919948679367
$wsEffluent.Range("B11").Value = 46.20916610716439
$wsEffluent.Range("B12").Value = 589.0053169068656
$wsEffluent.Range("B13").Value = 7024.509074669612
$wsEffluent.Range("B14").Value = 213.2732902092059
$wsEffluent.Range("B15").Value = 303.962645514275
$wsEffluent.Range("B16").Value = 301.410659170956
$wsEffluent.Range("B17").Value = 4946.454791319757
$wsEffluent.Range("B18").Value = 4801.899290309422
$wsEffluent.Range("B19").Value = 66.83164503716164
$wsEffluent.Range("B20").Value = 835.7166546854567
$wsEffluent.Range("B21").Value = 42.54580835492751
$wsEffluent.Range("B22").Value = 574.2341068625669
$wsEffluent.Range("B23").Value = 7084.607346499512
$wsEffluent.Range("B24").Value = 209.1045957988432
$wsEffluent.Range("B25").Value = 295.8378432015441
$wsEffluent.Range("B26").Value = 290.3098652538909
$wsEffluent.Range("B27").Value = 5018.514593746215
$wsEffluent.Range("B28").Value = 4946.454791319758
$wsEffluent.Range("B29").Value = 62.70281830017213
$wsEffluent.Range("B30").Value = 740.6919948679371
$wsEffluent.Range("B31").Value = 41.79184537105311
$wsEffluent.Range("B32").Value = 607.6013727739971
$wsEffluent.Range("B33").Value = 6692.694924610366
$wsEffluent.Range("B34").Value = 217.9980525548408
$wsEffluent.Range("B35").Value = 295.8378432015441
$wsEffluent.Range("B36").Value = 294.1929185923519
$wsEffluent.Range("B37").Value = 4946.454791319766
$wsEffluent.Range("B38").Value = 4659.596839528162
$wsEffluent.Range("B39").Value = 50.27758713229958
$wsEffluent.Range("B40").Value = 740.6919948679367
$wsEffluent.Range("B41").Value = 45.38455881536931
$wsEffluent.Range("B42").Value = 638.302376713069
$wsEffluent.Range("B43").Value = 8009.67552799563
$wsEffluent.Range("B44").Value = 221.8583735353482
$wsEffluent.Range("B45").Value = 323.648628513309
$wsEffluent.Range("B46").Value = 306.9835118461574
$wsEffluent.Range("B47").Value = 4946.454791319757
$wsEffluent.Range("B48").Value = 4946.454791319784
$wsEffluent.Range("B49").Value = 69.04319744715279
$wsEffluent.Range("B50").Value = 789.5667982073164
$wsEffluent.Range("B51").Value = 50.27758713229959
$wsEffluent.Range("B52").Value = 30.28014014042899
$wsEffluent.Range("B53").Value = 1061.288784583341
$wsEffluent.Range("B54").Value = 673.139019756272
$wsEffluent.Range("B55").Value = 13193.69723964875
$wsEffluent.Range("B56").Value = 30.95485279414541
$wsEffluent.Range("B57").Value = 386.102960640127
$wsEffluent.Range("B58").Value = 151.1830456348136
$wsEffluent.Range("B59").Value = 521.7652904644993
$wsEffluent.Range("B60").Value = 166.2676711344307
$wsEffluent.Range("B61").Value = 440.9951535542905
$wsEffluent.Range("B62").Value = 264.2862802638527
$wsEffluent.Range("B63").Value = 11174.43788457828
$wsEffluent.Range("B64").Value = 295.8378432015435
$wsEffluent.Range("B65").Value = 9589.54526217734
$wsEffluent.Range("B66").Value = 2.992437396175996
$wsEffluent.Range("B67").Value = 119.1646162768412
$wsEffluent.Range("B68").Value = 37.77594757766878
$wsEffluent.Range("B69").Value = 1583.386087511009
$wsEffluent.Range("B70").Value = 1.870552976500988
$wsEffluent.Range("B71").Value = 90.57472301947769
